$d = $word.ActiveDocument

# The three existing paragraphs under the "Sound" heading that describe
# images/sound atmosphere, the castle music and the ambient sounds are
# replaced by Melissa's new, much more detailed "Sound" write-up.
$firstPara = $d.Paragraphs.Item(58)
$lastPara  = $d.Paragraphs.Item(60)

$rng = $d.Range($firstPara.Range.Start, $lastPara.Range.End)

$p1 = "Das Spiel ist mit Musik unterlegt, die sich stetig passend zu den Orten, an denen sich der Spieler bewegt, verändert. "
$p2 = "Zu Beginn, wenn die eigene Spielfigur ihre Vergangenheit reflektiert, wird die Szene mit einer ruhigen Klaviermusik hinterlegt (Bsp. 1). Sobald sich der Spieler in Edinburgh aufhält, um auf ein Taxi zu warten, wechselt die Musik zu einer keltischen Melodie, die passend zum Ambiente der Stadt auf einem Dudelsack gespielt wird (Bsp. 2). Die Melodie ist zu diesem Zeitpunkt eher heiter und soll dem Spieler helfen, sich in die schottische Gegend, mit all ihren romantischen Gebäuden und vielen Geheimnissen hinein zu versetzen. "
$p3 = "Sobald man im Schloss angekommen ist, erzeugt die Musik eine bedrückende Atmosphäre, die die schaurige, leicht gruselige Stimmung in dem alten Schloss widerspiegelt (Bsp. 3). Nachdem der Tote aufgefunden wurde, verändert sich die Musik erneut, da Scarlett im Musikzimmer beginnt, Harfe zu spielen (Bsp. 4). Dass Scarlett dafür verantwortlich ist, weiß der Spieler zu diesem Zeitpunkt nicht. Sobald man das Musikzimmer betritt, macht der Spieler Scarlett darauf aufmerksam, dass ihre Lieder sehr traurig klingen. Daraufhin intoniert sie ein eher fröhlich gehaltenes Harfenstück. Dieses endet jedoch nach dem kurzen Vorspiel und Scarlett verlässt das Zimmer (Bsp. 5). "
$p4 = "Danach wechselt die Musik wieder zu einer eher schaurigen Hintergrundmusik (Bsp. 6). Sobald man Violas Zimmer betritt und ihre Spieluhr findet, wechselt die Musik erneut. Dieses Mal hört man eine gruselige Musik mit der Melodie einer Spieluhr (Bsp. 7). Dieser Musikwechsel soll das Spiel lebhafter erscheinen lassen. "
$p5 = "Gegen Ende des Spiels, wenn der Mordfall gelöst wurde und man mit dem Taxi durch die Highlands fährt, ändert sich die Musik wieder zu einer sehr traditionellen alten keltischen Weise (Bsp. 8)."

$rng.Text = $p1 + "`r" + $p2 + "`r" + $p3 + "`r" + $p4 + "`r" + $p5

# After the replacement, the five new paragraphs occupy slots 58-62.
# Paragraphs 58-61 (and the runs in 62) use Times New Roman for the East
# Asian / complex-script font slots; only the very last paragraph mark
# keeps the plain (no rFonts) formatting of the paragraph it replaced.
for ($i = 58; $i -le 61; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.NameFarEast = "Times New Roman"
    $p.Range.Font.NameBi = "Times New Roman"
}

$lastNew = $d.Paragraphs.Item(62)
$lastNew.Range.Font.NameFarEast = "Times New Roman"
$lastNew.Range.Font.NameBi = "Times New Roman"
# restore the plain formatting (sz/szCs only) on the trailing paragraph mark
$markRng = $d.Range($lastNew.Range.End - 1, $lastNew.Range.End)
$markRng.Font.NameFarEast = ""
$markRng.Font.NameBi = ""
